{"js": "// The document is a daily \"three-digit x one-digit multiplication\" worksheet:\n// one heading paragraph with the date, followed by a single table whose\n// non-blank rows each hold five \"NNN\u00d7N=\" problems (one run per cell).\n// The commit swaps in a new day's date plus a fresh batch of problems,\n// applied strictly in document order (a couple of the old problem texts\n// repeat verbatim but map to different replacements depending on position,\n// so we must walk the paragraphs in order rather than search-and-replace\n// by unique text).\n\n// New values, in the exact order they appear walking the document from\n// top to bottom: the date heading first, then every non-empty paragraph\n// (math-problem cell) row by row, left to right.\nconst newValues = [\n  \"2025-08-19 Tuesday\",\n  \"982\u00d73=\", \"133\u00d76=\", \"703\u00d75=\", \"466\u00d79=\", \"793\u00d74=\",\n  \"775\u00d77=\", \"202\u00d74=\", \"589\u00d74=\", \"791\u00d77=\", \"130\u00d74=\",\n  \"630\u00d77=\", \"483\u00d73=\", \"480\u00d74=\", \"571\u00d76=\", \"656\u00d78=\",\n  \"268\u00d76=\", \"624\u00d78=\", \"509\u00d79=\", \"741\u00d72=\", \"677\u00d74=\",\n  \"617\u00d73=\", \"267\u00d73=\", \"185\u00d77=\", \"970\u00d76=\", \"985\u00d76=\"\n];\n\n// context.document.body.paragraphs walks every paragraph in the body in\n// document order, including those nested inside table cells, so it lines\n// up exactly with the order the diff's replacements occur in.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet valueIndex = 0;\nfor (let i = 0; i < paragraphs.items.length && valueIndex < newValues.length; i++) {\n  const paragraph = paragraphs.items[i];\n  if (paragraph.text.trim() === \"\") {\n    continue; // blank worksheet rows (left for student work) are untouched\n  }\n  paragraph.insertText(newValues[valueIndex], \"Replace\");\n  valueIndex++;\n}\n\nawait context.sync();\n", "ps1": "# The document is a daily \"three-digit x one-digit multiplication\" worksheet:\n# one heading paragraph with the date, followed by a single table whose\n# non-blank rows each hold five \"NNN\u00d7N=\" problems (one run per cell).\n# The commit swaps in a new day's date plus a fresh batch of problems,\n# applied strictly in document order (a couple of the old problem texts\n# repeat verbatim but map to different replacements depending on position,\n# so we must walk the paragraphs in order rather than search-and-replace\n# by unique text).\n\n$d = $word.ActiveDocument\n\n# New values, in the exact order they appear walking the document from\n# top to bottom: the date heading first, then every non-empty paragraph\n# (math-problem cell) row by row, left to right.\n$newValues = @(\n    \"2025-08-19 Tuesday\",\n    \"982\u00d73=\", \"133\u00d76=\", \"703\u00d75=\", \"466\u00d79=\", \"793\u00d74=\",\n    \"775\u00d77=\", \"202\u00d74=\", \"589\u00d74=\", \"791\u00d77=\", \"130\u00d74=\",\n    \"630\u00d77=\", \"483\u00d73=\", \"480\u00d74=\", \"571\u00d76=\", \"656\u00d78=\",\n    \"268\u00d76=\", \"624\u00d78=\", \"509\u00d79=\", \"741\u00d72=\", \"677\u00d74=\",\n    \"617\u00d73=\", \"267\u00d73=\", \"185\u00d77=\", \"970\u00d76=\", \"985\u00d76=\"\n)\n\n# $d.Paragraphs walks every paragraph in the document in order, including\n# those nested inside table cells, so it lines up exactly with the order\n# the diff's replacements occur in. Each paragraph Range.Text carries a\n# trailing mark (paragraph mark \"\\r\", or a table cell's end-of-cell \"\\r\\a\")\n# that plain .Trim() does NOT strip (0x07 \"\\a\" isn't whitespace), so strip\n# those control characters explicitly to detect genuinely blank cells.\n$paragraphCount = $d.Paragraphs.Count\n$valueIndex = 0\nfor ($i = 1; $i -le $paragraphCount -and $valueIndex -lt $newValues.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n    $core = $r.Text.TrimEnd([char]13, [char]7)\n    if ($core.Length -eq 0) {\n        continue # blank worksheet rows (left for student work) are untouched\n    }\n    # Replace just the text, not the trailing paragraph/cell mark, so the\n    # run's formatting (font/size) and the paragraph structure survive.\n    $r.MoveEnd(1, -1) | Out-Null\n    $r.Text = $newValues[$valueIndex]\n    $valueIndex++\n}\n"}
